# Updated symbol list on Fri Dec 23 13:34:22 UTC 2022 with GitHub Actions
# Applies refreshed price/volume figures to the cryptos worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$CellRef,
        [string]$Text
    )
    $rng = $ws.Range($CellRef)
    # Force text storage so numeric-looking strings (e.g. "246.06") are not
    # silently coerced into numbers by Excel's type inference, then restore
    # the default "Normal" style so no stray style index is left behind.
    $rng.NumberFormat = "@"
    $rng.Value = $Text
    $rng.Style = "Normal"
}

# Price column (D) updates
Set-TextValue "D2"  "246.06"
Set-TextValue "D4"  "5.411"
Set-TextValue "D5"  "0.05870"
Set-TextValue "D6"  "3.382"
Set-TextValue "D7"  "6.360"
Set-TextValue "D8"  "0.8148"
Set-TextValue "D9"  "1.024"
Set-TextValue "D10" "0.01122"
Set-TextValue "D12" "0.07439"
Set-TextValue "D13" "0.03470"
Set-TextValue "D14" "0.03026"
Set-TextValue "D15" "4.196"
Set-TextValue "D16" "0.09405"
Set-TextValue "D17" "0.001594"
Set-TextValue "D18" "0.04810"
Set-TextValue "D19" "0.006062"
Set-TextValue "D20" "0.004114"
Set-TextValue "D21" "0.0009946"
Set-TextValue "D24" "2.223"
Set-TextValue "D25" "0.3244"
Set-TextValue "D27" "0.0001291"
Set-TextValue "D40" "0.03852"
Set-TextValue "D41" "0.006462"
Set-TextValue "D42" "0.1080"
Set-TextValue "D43" "0.002600"
Set-TextValue "D44" "0.006124"
Set-TextValue "D45" "0.00005625"
Set-TextValue "D47" "0.5303"
Set-TextValue "D48" "0.1423"

# Volume(1h) column (E) updates
$ws.Range("E10").Value = "9OneONEBestin24h"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"
